{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the \"GIS & Geospatial Analysis Consulting\" paragraph under the\n// Siege Analytics (PARTNER) role, right before its bullet list.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"GIS & Geospatial Analysis Consulting\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Paragraph \"GIS & Geospatial Analysis Consulting\" not found.');\n}\n\nconst newBullets = [\n  \"\\u2022 Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels\",\n  \"\\u2022 Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide\",\n  \"\\u2022 Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis\"\n];\n\n// Insert the three new bullet paragraphs, in order, immediately after the\n// \"GIS & Geospatial Analysis Consulting\" paragraph.\nlet anchor = target;\nfor (const bulletText of newBullets) {\n  anchor = anchor.insertParagraph(bulletText, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$bullet = [string][char]8226\n\n# Locate the \"GIS & Geospatial Analysis Consulting\" paragraph under the\n# Siege Analytics (PARTNER) role, right before its bullet list.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text.TrimEnd(\"`r`a\") -eq \"GIS & Geospatial Analysis Consulting\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Paragraph 'GIS & Geospatial Analysis Consulting' not found.\"\n}\n\n$newBullets = @(\n    ($bullet + \" Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels\"),\n    ($bullet + \" Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide\"),\n    ($bullet + \" Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis\")\n)\n\n# Insert the three new bullet paragraphs, in order, immediately after the\n# \"GIS & Geospatial Analysis Consulting\" paragraph.\n$insertIndex = $targetIndex\nforeach ($line in $newBullets) {\n    $d.Paragraphs($insertIndex).Range.InsertParagraphAfter()\n    $insertIndex = $insertIndex + 1\n    $d.Paragraphs($insertIndex).Range.Text = $line\n}\n"}
